$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 620, shifting the existing rows 620:661 down to 621:662.
$ws.Rows(620).Insert()

# Column A holds a date-like string (e.g. "2026/01/12") that must stay literal
# text, matching the sheet's existing inline-string date cells -- so force a
# text number format before writing the value (otherwise Excel auto-parses it
# into a date serial), then clear the temporary format back off the cell.
$ws.Range("A620").NumberFormat = "@"
$ws.Range("A620").Value = "2026/01/12"
$ws.Range("A620").ClearFormats()

$ws.Range("B620").Value = "月"
$ws.Range("C620").Value = 13
$ws.Range("D620").Value = 201
